$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# (matching the source data which stores them as text strings), then the
# temporary number-format style is reset back to Normal so no new style
# is left applied to the cell.
$textForceCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D12",
    "D15",
    "D18",
    "D22",
    "D23",
    "D24",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D45",
    "D48",
    "D50"
)
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "69.116.61"
$ws.Range("E2").Value = "  -1.89%  "

# Row 3
$ws.Range("D3").Value = "3.519.27"
$ws.Range("E3").Value = "  -2.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").Value = "582.78"
$ws.Range("E5").Value = "  -3.43%  "

# Row 6
$ws.Range("D6").Value = "192.88"
$ws.Range("E6").Value = "  -1.77%  "

# Row 7
$ws.Range("E7").Value = "  -3.49%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "0.204"
$ws.Range("E9").Value = "  -2.61%  "

# Row 10
$ws.Range("D10").Value = "0.618"
$ws.Range("E10").Value = "  -4.32%  "

# Row 11
$ws.Range("D11").Value = "52.10"
$ws.Range("E11").Value = "  -2.19%  "

# Row 12
$ws.Range("D12").Value = "0.0000286"
$ws.Range("E12").Value = "  -5.51%  "

# Row 13
$ws.Range("E13").Value = "  -4.54%  "

# Row 14
$ws.Range("D14").Value = "4.073.78"
$ws.Range("E14").Value = "  -3.02%  "

# Row 15
$ws.Range("D15").Value = "644.14"
$ws.Range("E15").Value = "  +7.35%  "

# Row 16
$ws.Range("D16").Value = "69.140.65"
$ws.Range("E16").Value = "  -2.08%  "

# Row 17
$ws.Range("D17").Value = "3.529.97"
$ws.Range("E17").Value = "  -2.17%  "

# Row 18
$ws.Range("D18").Value = "12.47"
$ws.Range("E18").Value = "  -3.72%  "

# Row 19
$ws.Range("E19").Value = "  -1.75%  "

# Row 20
$ws.Range("E20").Value = "  -4.28%  "

# Row 21
$ws.Range("E21").Value = "  -4.37%  "

# Row 22
$ws.Range("D22").Value = "17.94"
$ws.Range("E22").Value = "  -2.61%  "

# Row 23
$ws.Range("D23").Value = "5.40"
$ws.Range("E23").Value = "  +3.79%  "

# Row 24
$ws.Range("D24").Value = "101.40"
$ws.Range("E24").Value = "  -0.87%  "

# Row 25
$ws.Range("E25").Value = "  -6.07%  "

# Row 26
$ws.Range("E26").Value = "  -3.54%  "

# Row 27
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  -5.23%  "

# Row 28
$ws.Range("D28").Value = "9.39"
$ws.Range("E28").Value = "  -2.85%  "

# Row 29
$ws.Range("D29").Value = "32.75"
$ws.Range("E29").Value = "  -2.94%  "

# Row 30
$ws.Range("D30").Value = "6.72"
$ws.Range("E30").Value = "  -8.38%  "

# Row 31
$ws.Range("D31").Value = "4.07"
$ws.Range("E31").Value = "  -12.59%  "

# Row 32
$ws.Range("D32").Value = "11.64"
$ws.Range("E32").Value = "  -5.21%  "

# Row 33
$ws.Range("E33").Value = "  -6.88%  "

# Row 34
$ws.Range("D34").Value = "61.21"
$ws.Range("E34").Value = "  -3.54%  "

# Row 35
$ws.Range("D35").Value = "3.718.63"
$ws.Range("E35").Value = "  -5.37%  "

# Row 36
$ws.Range("E36").Value = "  -0.20%  "

# Row 37
$ws.Range("D37").Value = "0.0₃0795"
$ws.Range("E37").Value = "  -10.07%  "

# Row 38
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "3.60"
$ws.Range("E38").Value = "  +1.79%  "

# Row 39
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "504.77"
$ws.Range("E39").Value = "  -5.38%  "

# Row 40
$ws.Range("D40").Value = "2.92"
$ws.Range("E40").Value = "  -4.58%  "

# Row 41
$ws.Range("E41").Value = "  -5.57%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.133"
$ws.Range("E42").Value = "  -0.16%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "34.47"
$ws.Range("E43").Value = "  -6.45%  "

# Row 44
$ws.Range("D44").Value = "0.0442"
$ws.Range("E44").Value = "  -4.42%  "

# Row 45
$ws.Range("D45").Value = "3.35"
$ws.Range("E45").Value = "  -3.63%  "

# Row 46
$ws.Range("E46").Value = "  -1.58%  "

# Row 47
$ws.Range("E47").Value = "  -4.28%  "

# Row 48
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.39%  "

# Row 49
$ws.Range("E49").Value = "  -4.66%  "

# Row 50
$ws.Range("D50").Value = "2.68"
$ws.Range("E50").Value = "  +57.96%  "

# Row 51
$ws.Range("E51").Value = "  +1.25%  "

foreach ($ref in $textForceCells) {
    $ws.Range($ref).Style = "Normal"
}
